# Insert a new weekly record for "Berenjena" at row 352 in the
# "Vega Central Mapocho de Santiago" sheet, pushing the existing
# rows 352-366 down to 353-367 (and therefore extending the used
# range from A1:R366 to A1:R367).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 352..366 down by one, duplicating formatting from the
# row above (standard Excel Insert behaviour) so style s="2" on the
# date column carries down correctly.
$ws.Rows(352).Insert()

# Populate the newly inserted row with the new record's data.
$ws.Cells.Item(352, 1).Value = 9
$ws.Cells.Item(352, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(352, 3).Value = "Metropolitana"
$ws.Cells.Item(352, 4).Value = 45075
$ws.Cells.Item(352, 5).Value = 13
$ws.Cells.Item(352, 6).Value = 100112001
$ws.Cells.Item(352, 7).Value = "Berenjena"
$ws.Cells.Item(352, 8).Value = "Sin especificar"
$ws.Cells.Item(352, 9).Value = "Primera"
$ws.Cells.Item(352, 10).Value = 70
$ws.Cells.Item(352, 11).Value = 7000
$ws.Cells.Item(352, 12).Value = 8000
$ws.Cells.Item(352, 13).Value = 7500
$ws.Cells.Item(352, 14).Value = "$/caja 50 unidades"
$ws.Cells.Item(352, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(352, 16).Value = 150
$ws.Cells.Item(352, 17).Value = 50
$ws.Cells.Item(352, 18).Value = "Hortaliza"
